$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A8').Value = '23/11/2025'
$ws.Range('B8').Value = '16:18'
$ws.Range('C8').Value = '16:17'
$ws.Range('D8').Value = 'l''hopital'
$ws.Range('E8').Value = 'bossand'
$ws.Range('I8').Value = 'Non observable'
$ws.Range('K8').Value = 'ras'
$ws.Range('L8').Value = 'casas'
$ws.Range('M8').Value = 'transchool'
$ws.Range('W8').Value = 'T5'
$ws.Range('X8').Value = 'beau'
$ws.Range('Y8').Value = 22635
$ws.Range('Z8').Value = 'Conforme'
$ws.Range('AA8').Value = 'Conforme'
$ws.Range('AB8').Value = 'Conforme'
$ws.Range('AC8').Value = 'Conforme'
$ws.Range('AD8').Value = 'Conforme'
$ws.Range('AE8').Value = 'Conforme'
$ws.Range('AF8').Value = 'Propre'
$ws.Range('AG8').Value = 'ras'
$ws.Range('AH8').Value = 'Conforme'
$ws.Range('AI8').Value = 'Conforme'
$ws.Range('AJ8').Value = 'Conforme'
$ws.Range('AK8').Value = 'Propre'
$ws.Range('AL8').Value = 'Propre'
$ws.Range('AM8').Value = 'Propre'
$ws.Range('AN8').Value = 'Propre'
$ws.Range('AO8').Value = 'ras'
$ws.Range('AP8').Value = 32
$ws.Range('AQ8').Value = 2
$ws.Range('AR8').Value = 'lebon'

$ws.Range('A9').Value = '23/11/2025'
$ws.Range('B9').Value = '16:23'
$ws.Range('C9').Value = '16:22'
$ws.Range('D9').Value = 'poli'
$ws.Range('E9').Value = 'bossand'
$ws.Range('I9').Value = 'Non observable'
$ws.Range('K9').Value = 'ras'
$ws.Range('L9').Value = 'casc'
$ws.Range('O9').Value = 'Sc'
$ws.Range('R9').Value = 'PT01'
$ws.Range('X9').Value = 'beau'
$ws.Range('Y9').Value = 20325
$ws.Range('Z9').Value = 'Conforme'
$ws.Range('AA9').Value = 'Conforme'
$ws.Range('AB9').Value = 'Conforme'
$ws.Range('AC9').Value = 'Conforme'
$ws.Range('AD9').Value = 'Conforme'
$ws.Range('AE9').Value = 'Conforme'
$ws.Range('AF9').Value = 'Propre'
$ws.Range('AG9').Value = 'ras'
$ws.Range('AH9').Value = 'Conforme'
$ws.Range('AI9').Value = 'Non observable'
$ws.Range('AJ9').Value = 'Conforme'
$ws.Range('AK9').Value = 'Propre'
$ws.Range('AL9').Value = 'Propre'
$ws.Range('AM9').Value = 'Propre'
$ws.Range('AN9').Value = 'Propre'
$ws.Range('AO9').Value = 'ras'
$ws.Range('AP9').Value = 10
$ws.Range('AQ9').Value = 2
$ws.Range('AR9').Value = 'lebon'

$ws.Range('A10').Value = '24/11/2025'
$ws.Range('B10').Value = '16:30'
$ws.Range('C10').Value = '16:26'
$ws.Range('D10').Value = 'momerstroff'
$ws.Range('E10').Value = 'bossand'
$ws.Range('I10').Value = 'Non observable'
$ws.Range('K10').Value = 'ras'
$ws.Range('L10').Value = 'rgeFluo57'
$ws.Range('N10').Value = 'Lr'
$ws.Range('S10').Value = "'21"
$ws.Range('S10').Style = "Normal"
$ws.Range('X10').Value = 'beau'
$ws.Range('Y10').Value = 22320
$ws.Range('Z10').Value = 'Conforme'
$ws.Range('AA10').Value = 'Conforme'
$ws.Range('AB10').Value = 'Conforme'
$ws.Range('AC10').Value = 'Conforme'
$ws.Range('AD10').Value = 'Conforme'
$ws.Range('AE10').Value = 'Conforme'
$ws.Range('AF10').Value = 'Propre'
$ws.Range('AG10').Value = 'ras'
$ws.Range('AH10').Value = 'Conforme'
$ws.Range('AI10').Value = 'Non observable'
$ws.Range('AJ10').Value = 'Conforme'
$ws.Range('AK10').Value = 'Propre'
$ws.Range('AL10').Value = 'Propre'
$ws.Range('AM10').Value = 'Propre'
$ws.Range('AN10').Value = 'Propre'
$ws.Range('AO10').Value = 'ras'
$ws.Range('AP10').Value = 20
$ws.Range('AQ10').Value = 1
$ws.Range('AR10').Value = 'lebon'

$ws.Range('A11').Value = '24/11/2025'
$ws.Range('B11').Value = '18:10'
$ws.Range('C11').Value = '18:08'
$ws.Range('D11').Value = 'Forbach'
$ws.Range('E11').Value = 'Bangoura'
$ws.Range('I11').Value = 'Non observable'
$ws.Range('K11').Value = 'ras'
$ws.Range('L11').Value = 'forbus'
$ws.Range('X11').Value = 'beau'
$ws.Range('Y11').Value = 22325
$ws.Range('Z11').Value = 'Conforme'
$ws.Range('AA11').Value = 'Conforme'
$ws.Range('AB11').Value = 'Conforme'
$ws.Range('AC11').Value = 'Conforme'
$ws.Range('AD11').Value = 'Conforme'
$ws.Range('AE11').Value = 'Conforme'
$ws.Range('AF11').Value = 'Propre'
$ws.Range('AG11').Value = 'ras'
$ws.Range('AH11').Value = 'Conforme'
$ws.Range('AI11').Value = 'Conforme'
$ws.Range('AJ11').Value = 'Conforme'
$ws.Range('AK11').Value = 'Propre'
$ws.Range('AL11').Value = 'Propre'
$ws.Range('AM11').Value = 'Propre'
$ws.Range('AN11').Value = 'Propre'
$ws.Range('AO11').Value = 'ras'
$ws.Range('AP11').Value = 10
$ws.Range('AQ11').Value = 1
$ws.Range('AR11').Value = 'lebon'

$ws.Range('A12').Value = '23/11/2025'
$ws.Range('B12').Value = '18:44'
$ws.Range('C12').Value = '18:43'
$ws.Range('D12').Value = 'Rosbruck'
$ws.Range('E12').Value = 'chahid'
$ws.Range('I12').Value = 'Non observable'
$ws.Range('K12').Value = 'ras'
$ws.Range('L12').Value = 'rgeFluo57'
$ws.Range('N12').Value = 'Lr'
$ws.Range('S12').Value = 'MS'
$ws.Range('X12').Value = 'beau'
$ws.Range('Y12').Value = 19965
$ws.Range('Z12').Value = 'Conforme'
$ws.Range('AA12').Value = 'Conforme'
$ws.Range('AB12').Value = 'Conforme'
$ws.Range('AC12').Value = 'Conforme'
$ws.Range('AD12').Value = 'Conforme'
$ws.Range('AE12').Value = 'Conforme'
$ws.Range('AF12').Value = 'Propre'
$ws.Range('AG12').Value = 'ras'
$ws.Range('AH12').Value = 'Conforme'
$ws.Range('AI12').Value = 'Conforme'
$ws.Range('AJ12').Value = 'Conforme'
$ws.Range('AK12').Value = 'Propre'
$ws.Range('AL12').Value = 'Propre'
$ws.Range('AM12').Value = 'Propre'
$ws.Range('AN12').Value = 'Propre'
$ws.Range('AO12').Value = "'"
$ws.Range('AO12').Style = "Normal"
$ws.Range('AP12').Value = 32
$ws.Range('AQ12').Value = 1
$ws.Range('AR12').Value = 'lebon'

$ws.Range('A13').Value = '23/11/2025'
$ws.Range('B13').Value = '19:10'
$ws.Range('C13').Value = '19:09'
$ws.Range('D13').Value = 'poli'
$ws.Range('E13').Value = 'Bangoura'
$ws.Range('I13').Value = 'Non observable'
$ws.Range('K13').Value = 'ras'
$ws.Range('L13').Value = 'rgeFluo57'
$ws.Range('N13').Value = 'Sa'
$ws.Range('T13').Value = 'SA7'
$ws.Range('X13').Value = 'beau'
$ws.Range('Y13').Value = 22352
$ws.Range('Z13').Value = 'Conforme'
$ws.Range('AA13').Value = 'Conforme'
$ws.Range('AB13').Value = 'Conforme'
$ws.Range('AC13').Value = 'Conforme'
$ws.Range('AD13').Value = 'Conforme'
$ws.Range('AE13').Value = 'Conforme'
$ws.Range('AF13').Value = 'Moyen'
$ws.Range('AG13').Value = 'ras'
$ws.Range('AH13').Value = 'Conforme'
$ws.Range('AI13').Value = 'Conforme'
$ws.Range('AJ13').Value = 'Conforme'
$ws.Range('AK13').Value = 'Propre'
$ws.Range('AL13').Value = 'Propre'
$ws.Range('AM13').Value = 'Propre'
$ws.Range('AN13').Value = 'Propre'
$ws.Range('AO13').Value = 'ras'
$ws.Range('AP13').Value = 10
$ws.Range('AQ13').Value = 2
$ws.Range('AR13').Value = 'lebon'

